$wb = $excel.ActiveWorkbook

# Sheet "1er Parcial" - row 7 (Morales Vallejo Jorge Luis / 2APV)
$ws1 = $wb.Worksheets.Item("1er Parcial")
$ws1.Range("E7").Value = 19
$ws1.Range("F7").Value = 6
$ws1.Range("G7").Value = 76
$ws1.Range("H7").Value = 24
$ws1.Range("I7").Value = 7.8
$ws1.Range("J7").Value = 6
$ws1.Range("K7").Value = 24

# Sheet "2o Parcial" - row 7
$ws2 = $wb.Worksheets.Item("2o Parcial")
$ws2.Range("E7").Value = 14
$ws2.Range("F7").Value = 11
$ws2.Range("G7").Value = 56
$ws2.Range("H7").Value = 44
$ws2.Range("I7").Value = 7
$ws2.Range("J7").Value = 11
$ws2.Range("K7").Value = 44

# Sheet "3er Parcial" - row 7
$ws3 = $wb.Worksheets.Item("3er Parcial")
$ws3.Range("E7").Value = 19
$ws3.Range("F7").Value = 6
$ws3.Range("G7").Value = 76
$ws3.Range("H7").Value = 24
$ws3.Range("I7").Value = 7.6
$ws3.Range("J7").Value = 6
$ws3.Range("K7").Value = 24
